# Add data for 2024-07-16: refresh the 2024 (column K) violent-crime
# counts across the citywide summary, the by-neighborhood rollup, and
# every individual neighborhood sheet. A handful of rows also pick up
# small 2021-2023 (H/I/J) corrections that shipped in the same refresh.
$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 4269
$ws.Range('K3').Value = 4343
$ws.Range('H4').Value = 1738
$ws.Range('I4').Value = 1797
$ws.Range('J4').Value = 1824
$ws.Range('K4').Value = 881
$ws.Range('K5').Value = 320
$ws.Range('K6').Value = 4873
$ws.Range('H7').Value = 26051
$ws.Range('I7').Value = 26252
$ws.Range('J7').Value = 29294
$ws.Range('K7').Value = 14686

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 277
$ws.Range('K3').Value = 294
$ws.Range('K4').Value = 57
$ws.Range('K6').Value = 330
$ws.Range('K7').Value = 985

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 110
$ws.Range('K3').Value = 108
$ws.Range('K7').Value = 312

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 228
$ws.Range('K6').Value = 177
$ws.Range('K7').Value = 612

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K5').Value = 10
$ws.Range('K7').Value = 254

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 138
$ws.Range('K5').Value = 25
$ws.Range('K7').Value = 500

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 101
$ws.Range('K6').Value = 133
$ws.Range('K7').Value = 336

# Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('K2').Value = 13
$ws.Range('K6').Value = 15
$ws.Range('K7').Value = 52

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 127
$ws.Range('K7').Value = 430
$ws.Range('K8').Value = 985
$ws.Range('K11').Value = 290
$ws.Range('K13').Value = 16
$ws.Range('K16').Value = 47
$ws.Range('K19').Value = 443
$ws.Range('K20').Value = 333
$ws.Range('K23').Value = 151
$ws.Range('K27').Value = 140
$ws.Range('K29').Value = 780
$ws.Range('K30').Value = 52
$ws.Range('K31').Value = 160
$ws.Range('K33').Value = 612
$ws.Range('K36').Value = 186
$ws.Range('K37').Value = 500
$ws.Range('K42').Value = 538
$ws.Range('K44').Value = 133
$ws.Range('K47').Value = 88
$ws.Range('K48').Value = 190
$ws.Range('K49').Value = 86
$ws.Range('K51').Value = 184
$ws.Range('K52').Value = 398
$ws.Range('K54').Value = 272
$ws.Range('K55').Value = 163
$ws.Range('H63').Value = 289
$ws.Range('I63').Value = 214
$ws.Range('J63').Value = 109
$ws.Range('K63').Value = 53
$ws.Range('K64').Value = 90
$ws.Range('K65').Value = 336
$ws.Range('K67').Value = 568
$ws.Range('K73').Value = 133
$ws.Range('K77').Value = 106
$ws.Range('K82').Value = 15
$ws.Range('K83').Value = 312
$ws.Range('K85').Value = 661
$ws.Range('K86').Value = 99
$ws.Range('K91').Value = 162
$ws.Range('K94').Value = 183
$ws.Range('K95').Value = 254
$ws.Range('K97').Value = 123
$ws.Range('H101').Value = 26051
$ws.Range('I101').Value = 26252
$ws.Range('J101').Value = 29294
$ws.Range('K101').Value = 14686

# Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 56
$ws.Range('K6').Value = 54
$ws.Range('K7').Value = 160

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K6').Value = 166
$ws.Range('K7').Value = 568

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K2').Value = 15
$ws.Range('K7').Value = 86

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K4').Value = 13
$ws.Range('K6').Value = 137
$ws.Range('K7').Value = 272

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 224
$ws.Range('K3').Value = 276
$ws.Range('K4').Value = 41
$ws.Range('K6').Value = 216
$ws.Range('K7').Value = 780

# Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K6').Value = 93
$ws.Range('K7').Value = 190

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 142
$ws.Range('K7').Value = 443

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K4').Value = 8
$ws.Range('K7').Value = 133

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 171
$ws.Range('K4').Value = 22
$ws.Range('K6').Value = 194
$ws.Range('K7').Value = 538

# Boystown
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('K5').Value = 8
$ws.Range('K6').Value = 16

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K3').Value = 45
$ws.Range('K7').Value = 163

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 43
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 151

# Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 76
$ws.Range('K6').Value = 39
$ws.Range('K7').Value = 162

# Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 90

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 116
$ws.Range('K3').Value = 104
$ws.Range('K7').Value = 333

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K6').Value = 44
$ws.Range('K7').Value = 186

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K6').Value = 105
$ws.Range('K7').Value = 430

# West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K4').Value = 18
$ws.Range('K7').Value = 183

# Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K2').Value = 29
$ws.Range('K7').Value = 88

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 94
$ws.Range('K3').Value = 73
$ws.Range('K7').Value = 290

# Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K3').Value = 32
$ws.Range('K7').Value = 133

# Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K6').Value = 42
$ws.Range('K7').Value = 127

# West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K3').Value = 22
$ws.Range('K7').Value = 123

# Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K6').Value = 54
$ws.Range('K7').Value = 140

# Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 38
$ws.Range('K7').Value = 99

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 52
$ws.Range('K7').Value = 184

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 233
$ws.Range('K3').Value = 221
$ws.Range('K6').Value = 155
$ws.Range('K7').Value = 661

# Sheffield & DePaul
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('K3').Value = 3
$ws.Range('K6').Value = 15

# Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K6').Value = 15
$ws.Range('K7').Value = 106

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 103
$ws.Range('K3').Value = 105
$ws.Range('K7').Value = 398

# Bucktown
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K3').Value = 5
$ws.Range('K7').Value = 47
